$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 399.5
$ws.Range("I12").Value = 399.5
$ws.Range("K12").Value = 399.5
$ws.Range("M12").Value = -229.5
$ws.Range("H61").Value = 120.28571
$ws.Range("I61").Value = 120.28571
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 360.85713
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -188.85713
$ws.Range("N61").ClearContents()
$ws.Range("H92").Value = 1188.0834
$ws.Range("I92").Value = 1041.1875
$ws.Range("J92").Value = 1481.875
$ws.Range("K92").Value = 1041.1875
$ws.Range("L92").Value = 1481.875
$ws.Range("M92").Value = 206.8125
$ws.Range("N92").Value = -3977.875
$ws.Range("H137").Value = 1465.6061
$ws.Range("I137").Value = 1056.9375
$ws.Range("K137").Value = 3170.8125
$ws.Range("M137").Value = -620.8125
$ws.Range("H138").Value = 507101.8
$ws.Range("I138").Value = 1150.6666
$ws.Range("J138").Value = 641030.0600000001
$ws.Range("K138").Value = 3451.9998
$ws.Range("L138").Value = 1923090.18
$ws.Range("M138").Value = 1688.0002
$ws.Range("N138").Value = -1933370.18
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1301.2222
$ws.Range("I2").Value = 1242.2
$ws.Range("J2").Value = 1375
$ws.Range("K2").Value = 1242.2
$ws.Range("L2").Value = 1375
$ws.Range("M2").Value = -1129.2
$ws.Range("N2").Value = -1601
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H32").Value = 4886.0303
$ws.Range("I32").Value = 4774.65
$ws.Range("J32").Value = 5999.8335
$ws.Range("K32").Value = 4774.65
$ws.Range("L32").Value = 5999.8335
$ws.Range("M32").Value = -4487.65
$ws.Range("N32").Value = -6573.8335
$ws.Range("H40").Value = 18000
$ws.Range("J40").Value = 18000
$ws.Range("L40").Value = 18000
$ws.Range("N40").Value = -18352
$ws.Range("H45").Value = 1809.8
$ws.Range("I45").Value = 1782.0769
$ws.Range("K45").Value = 1782.0769
$ws.Range("M45").Value = -1405.0769
$ws.Range("H61").Value = 55556850
$ws.Range("I61").Value = 83334296
$ws.Range("J61").Value = 1950.3334
$ws.Range("K61").Value = 83334296
$ws.Range("L61").Value = 1950.3334
$ws.Range("M61").Value = -83334084
$ws.Range("N61").Value = -2374.3334
$ws.Range("H116").Value = 1301.2222
$ws.Range("I116").Value = 1242.2
$ws.Range("J116").Value = 1375
$ws.Range("K116").Value = 1242.2
$ws.Range("L116").Value = 1375
$ws.Range("M116").Value = 1051.8
$ws.Range("N116").Value = -5963
$ws.Range("H136").Value = 55556850
$ws.Range("I136").Value = 83334296
$ws.Range("J136").Value = 1950.3334
$ws.Range("K136").Value = 250002888
$ws.Range("L136").Value = 5851.0002
$ws.Range("M136").Value = -250000338
$ws.Range("N136").Value = -10951.0002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1301.2222
$ws.Range("I3").Value = 1242.2
$ws.Range("J3").Value = 1375
$ws.Range("K3").Value = 1242.2
$ws.Range("L3").Value = 1375
$ws.Range("M3").Value = -1128.2
$ws.Range("N3").Value = -1603
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H37").Value = 3000
$ws.Range("I37").Value = 666.6667
$ws.Range("K37").Value = 666.6667
$ws.Range("M37").Value = -529.6667
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H127").Value = 42596
$ws.Range("J127").Value = 42596
$ws.Range("L127").Value = 42596
$ws.Range("N127").Value = -52516

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1523.3448
$ws.Range("I58").Value = 1269.7778
$ws.Range("J58").Value = 1938.2727
$ws.Range("K58").Value = 1269.7778
$ws.Range("L58").Value = 1938.2727
$ws.Range("M58").Value = -1066.7778
$ws.Range("N58").Value = -2344.2727
$ws.Range("H74").Value = 33000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 33000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H108").Value = 32786
$ws.Range("J108").Value = 32786
$ws.Range("L108").Value = 32786
$ws.Range("N108").Value = -40466
$ws.Range("H136").Value = 1523.3448
$ws.Range("I136").Value = 1269.7778
$ws.Range("J136").Value = 1938.2727
$ws.Range("K136").Value = 3809.3334
$ws.Range("L136").Value = 5814.8181
$ws.Range("M136").Value = -1259.3334
$ws.Range("N136").Value = -10914.8181
$ws.Range("H141").Value = 372130.6
$ws.Range("J141").Value = 372130.6
$ws.Range("L141").Value = 372130.6
$ws.Range("N141").Value = -382490.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 449.25
$ws.Range("I5").Value = 369.82352
$ws.Range("J5").Value = 642.1429000000001
$ws.Range("K5").Value = 1109.47056
$ws.Range("L5").Value = 1926.4287
$ws.Range("M5").Value = -997.47056
$ws.Range("N5").Value = -2150.4287
$ws.Range("H70").Value = 4829.6875
$ws.Range("I70").Value = 1895.2
$ws.Range("K70").Value = 5685.6
$ws.Range("M70").Value = -5370.6
$ws.Range("H73").Value = 4829.6875
$ws.Range("I73").Value = 1895.2
$ws.Range("K73").Value = 5685.6
$ws.Range("M73").Value = -4593.6
$ws.Range("H105").Value = 89486.586
$ws.Range("J105").Value = 89486.586
$ws.Range("L105").Value = 268459.758
$ws.Range("N105").Value = -273701.758
$ws.Range("H117").Value = 768.7895
$ws.Range("I117").Value = 500.5
$ws.Range("J117").Value = 1066.8889
$ws.Range("K117").Value = 1501.5
$ws.Range("L117").Value = 3200.6667
$ws.Range("M117").Value = 1940.5
$ws.Range("N117").Value = -10084.6667
$ws.Range("H131").Value = 17858058
$ws.Range("I131").Value = 71428984
$ws.Range("J131").Value = 1083.2142
$ws.Range("K131").Value = 214286952
$ws.Range("L131").Value = 3249.6426
$ws.Range("M131").Value = -214281912
$ws.Range("N131").Value = -13329.6426
$ws.Range("H135").Value = 449.25
$ws.Range("I135").Value = 369.82352
$ws.Range("J135").Value = 642.1429000000001
$ws.Range("K135").Value = 3328.41168
$ws.Range("L135").Value = 5779.2861
$ws.Range("M135").Value = -793.4116799999997
$ws.Range("N135").Value = -10849.2861

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3504.25
$ws.Range("I43").Value = 3005.6667
$ws.Range("K43").Value = 3005.6667
$ws.Range("M43").Value = -2854.6667
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 10000
$ws.Range("N48").Value = -10970
$ws.Range("H49").Value = 23000
$ws.Range("J49").Value = 23000
$ws.Range("L49").Value = 23000
$ws.Range("N49").Value = -23368
$ws.Range("H134").Value = 26754.555
$ws.Range("J134").Value = 26754.555
$ws.Range("L134").Value = 80263.66500000001
$ws.Range("N134").Value = -85333.66500000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3328.6428
$ws.Range("I40").Value = 1971.7778
$ws.Range("J40").Value = 5771
$ws.Range("K40").Value = 1971.7778
$ws.Range("L40").Value = 5771
$ws.Range("M40").Value = -1835.7778
$ws.Range("N40").Value = -6043
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H135").Value = 40921.8
$ws.Range("J135").Value = 40921.8
$ws.Range("L135").Value = 40921.8
$ws.Range("N135").Value = -51061.8
$ws.Range("H136").Value = 1567
$ws.Range("I136").Value = 1223.5714
$ws.Range("J136").Value = 2368.3333
$ws.Range("K136").Value = 3670.7142
$ws.Range("L136").Value = 7104.999899999999
$ws.Range("M136").Value = -1120.7142
$ws.Range("N136").Value = -12204.9999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 111115350
$ws.Range("I126").Value = 125003520
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 375010560
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -375008090
$ws.Range("N126").Value = -34940
